# Append a new row of data (row 7) to Sheet1, matching the existing
# columns: MIGRATION DATE / FINANCIAL INSTITUTION NAME / ENTITY ID / ADDRESS
# style rows above (A: Chennai, B: Tamil Nadu, C: 600117 pin code, D: phone).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "Chennai"
$ws.Range("B7").Value = "Tamil Nadu"
# Pin code and phone number are textual identifiers (leading context implies
# text, like the other inline-string columns), so force text entry with a
# leading apostrophe to avoid Excel auto-converting them to numbers.
$ws.Range("C7").Value = "'600117"
$ws.Range("D7").Value = "'9911991100"
